$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(149, 2).Value = 48654
$ws.Cells.Item(149, 5).Value = 38.26
$ws.Cells.Item(149, 6).Value = -1
$ws.Cells.Item(149, 7).Value = -32.02
$ws.Cells.Item(150, 2).Value = 63902
$ws.Cells.Item(150, 5).Value = 34.04
$ws.Cells.Item(150, 6).Value = 2
$ws.Cells.Item(150, 7).Value = 64.04000000000001
$ws.Cells.Item(161, 2).Value = 57756
$ws.Cells.Item(161, 5).Value = 79.37
$ws.Cells.Item(161, 6).Value = -100
$ws.Cells.Item(161, 7).Value = -6644
$ws.Cells.Item(162, 2).Value = 53925
$ws.Cells.Item(162, 6).Value = 1
$ws.Cells.Item(162, 7).Value = 66.44
$ws.Cells.Item(163, 2).Value = 64350
$ws.Cells.Item(163, 5).Value = 70.63
$ws.Cells.Item(163, 6).Value = 101
$ws.Cells.Item(163, 7).Value = 6710.44
$ws.Cells.Item(183, 2).Value = 57552
$ws.Cells.Item(183, 5).Value = 136.86
$ws.Cells.Item(183, 6).Value = -5
$ws.Cells.Item(183, 7).Value = -603.45
$ws.Cells.Item(184, 2).Value = 64329
$ws.Cells.Item(184, 5).Value = 128.32
$ws.Cells.Item(184, 6).Value = 6
$ws.Cells.Item(184, 7).Value = 724.14
$ws.Cells.Item(264, 2).Value = 48719
$ws.Cells.Item(264, 5).Value = 353.35
$ws.Cells.Item(264, 6).Value = -81
$ws.Cells.Item(264, 7).Value = -23955.75
$ws.Cells.Item(265, 2).Value = 64979
$ws.Cells.Item(265, 5).Value = 314.41
$ws.Cells.Item(265, 6).Value = 82
$ws.Cells.Item(265, 7).Value = 24251.5
$ws.Cells.Item(279, 2).Value = 64973
$ws.Cells.Item(279, 5).Value = 35.4
$ws.Cells.Item(279, 6).Value = 150
$ws.Cells.Item(279, 7).Value = 4995
$ws.Cells.Item(280, 2).Value = 48706
$ws.Cells.Item(280, 5).Value = 39.8
$ws.Cells.Item(280, 6).Value = -144
$ws.Cells.Item(280, 7).Value = -4795.2
$ws.Cells.Item(316, 2).Value = 61610
$ws.Cells.Item(316, 4).Value = 102.71
$ws.Cells.Item(316, 5).Value = 122.71
$ws.Cells.Item(316, 6).Value = -58
$ws.Cells.Item(316, 7).Value = -5957.18
$ws.Cells.Item(318, 2).Value = 57077
$ws.Cells.Item(318, 4).Value = 93.08
$ws.Cells.Item(318, 5).Value = 111.2
$ws.Cells.Item(318, 6).Value = 1
$ws.Cells.Item(318, 7).Value = 93.08
$ws.Cells.Item(346, 2).Value = 55373
$ws.Cells.Item(346, 5).Value = 163.62
$ws.Cells.Item(346, 6).Value = -94
$ws.Cells.Item(346, 7).Value = -13562.32
$ws.Cells.Item(347, 2).Value = 63520
$ws.Cells.Item(347, 5).Value = 153.4
$ws.Cells.Item(347, 6).Value = 97
$ws.Cells.Item(347, 7).Value = 13995.16
$ws.Cells.Item(351, 2).Value = 63571
$ws.Cells.Item(351, 5).Value = 152.53
$ws.Cells.Item(351, 6).Value = 29
$ws.Cells.Item(351, 7).Value = 4160.92
$ws.Cells.Item(352, 2).Value = 57802
$ws.Cells.Item(352, 5).Value = 162.71
$ws.Cells.Item(352, 6).Value = -79
$ws.Cells.Item(352, 7).Value = -11334.92
$ws.Cells.Item(355, 2).Value = 55356
$ws.Cells.Item(355, 5).Value = 54.04
$ws.Cells.Item(355, 6).Value = -158
$ws.Cells.Item(355, 7).Value = -7527.12
$ws.Cells.Item(356, 2).Value = 63510
$ws.Cells.Item(356, 5).Value = 50.66
$ws.Cells.Item(356, 6).Value = 167
$ws.Cells.Item(356, 7).Value = 7955.88
$ws.Cells.Item(400, 2).Value = 62933
$ws.Cells.Item(400, 6).Value = 146
$ws.Cells.Item(400, 7).Value = 8632.98
$ws.Cells.Item(401, 2).Value = 57835
$ws.Cells.Item(401, 6).Value = 1
$ws.Cells.Item(401, 7).Value = 59.13
$ws.Cells.Item(431, 2).Value = 63102
$ws.Cells.Item(431, 3).Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Cells.Item(431, 6).Value = 36
$ws.Cells.Item(431, 7).Value = 2140.92
$ws.Cells.Item(432, 2).Value = 53082
$ws.Cells.Item(432, 3).Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Cells.Item(432, 6).Value = 1
$ws.Cells.Item(432, 7).Value = 59.47
$ws.Cells.Item(579, 2).Value = 65069
$ws.Cells.Item(579, 5).Value = 14.3
$ws.Cells.Item(579, 6).Value = 172
$ws.Cells.Item(579, 7).Value = 2313.4
$ws.Cells.Item(580, 2).Value = 53757
$ws.Cells.Item(580, 5).Value = 16.08
$ws.Cells.Item(580, 6).Value = -159
$ws.Cells.Item(580, 7).Value = -2138.55
$ws.Cells.Item(583, 2).Value = 65066
$ws.Cells.Item(583, 5).Value = 13.61
$ws.Cells.Item(583, 6).Value = 313
$ws.Cells.Item(583, 7).Value = 4009.53
$ws.Cells.Item(584, 2).Value = 53263
$ws.Cells.Item(584, 5).Value = 15.29
$ws.Cells.Item(584, 6).Value = -309
$ws.Cells.Item(584, 7).Value = -3958.29
$ws.Cells.Item(586, 2).Value = 45695
$ws.Cells.Item(586, 5).Value = 23.58
$ws.Cells.Item(586, 6).Value = -36
$ws.Cells.Item(586, 7).Value = -710.28
$ws.Cells.Item(587, 2).Value = 64915
$ws.Cells.Item(587, 5).Value = 20.98
$ws.Cells.Item(587, 6).Value = 40
$ws.Cells.Item(587, 7).Value = 789.2
$ws.Cells.Item(599, 2).Value = 64925
$ws.Cells.Item(599, 5).Value = 13.97
$ws.Cells.Item(599, 6).Value = 302
$ws.Cells.Item(599, 7).Value = 3971.3
$ws.Cells.Item(600, 2).Value = 45709
$ws.Cells.Item(600, 5).Value = 15.69
$ws.Cells.Item(600, 6).Value = -300
$ws.Cells.Item(600, 7).Value = -3945
$ws.Cells.Item(687, 2).Value = 53319
$ws.Cells.Item(687, 5).Value = 310.64
$ws.Cells.Item(687, 6).Value = -6
$ws.Cells.Item(687, 7).Value = -1643.52
$ws.Cells.Item(688, 2).Value = 64810
$ws.Cells.Item(688, 5).Value = 291.22
$ws.Cells.Item(688, 6).Value = 7
$ws.Cells.Item(688, 7).Value = 1917.44
